$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 144.55556
Write-Host ("H4 now: " + $ws.Range("H4").Value)
